# Update the "cryptos" worksheet with refreshed price/volume figures,
# and fix the swapped EnergySwap/Algorand rows (50 and 51).
#
# Numeric-looking D-column values (e.g. "212.60", "0.535") are written via
# a temporary Text number format so Excel keeps the exact original string
# (including trailing zeros) instead of silently coercing it to a number;
# the format/style is then reset to Normal so no stray formatting remains.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Cells.Item(2, 4).Value = '27.549.41'
$ws.Cells.Item(2, 5).Value = '  -0.25%  '

# --- Row 3: Ethereum ---
$ws.Cells.Item(3, 4).Value = '1.647.79'
$ws.Cells.Item(3, 5).Value = '  -0.69%  '

# --- Row 4: TetherUSD ---
$ws.Cells.Item(4, 5).Value = '  +0.09%  '

# --- Row 5: BNB ---
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '212.60'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.18%  '

# --- Row 6: XRP ---
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.535'
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +4.90%  '

# --- Row 7: USDC ---
$ws.Cells.Item(7, 5).Value = '  +0.09%  '

# --- Row 8: Solana ---
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '23.56'
$ws.Cells.Item(8, 4).Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -2.54%  '

# --- Row 9: Cardano ---
$ws.Cells.Item(9, 5).Value = '  -1.67%  '

# --- Row 10: Dogecoin ---
$ws.Cells.Item(10, 5).Value = '  -1.34%  '

# --- Row 11: TRON ---
$ws.Cells.Item(11, 5).Value = '  +1.32%  '

# --- Row 12: Wrapped liquid staked Ether 2.0 ---
$ws.Cells.Item(12, 4).Value = '1.881.41'
$ws.Cells.Item(12, 5).Value = '  -0.64%  '

# --- Row 13: Wrapped Ether ---
$ws.Cells.Item(13, 4).Value = '1.644.53'
$ws.Cells.Item(13, 5).Value = '  +0.01%  '

# --- Row 15: Polkadot ---
$ws.Cells.Item(15, 5).Value = '  -2.48%  '

# --- Row 16: Litecoin ---
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.57'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -2.08%  '

# --- Row 17: Wrapped BTC ---
$ws.Cells.Item(17, 4).Value = '27.528.56'
$ws.Cells.Item(17, 5).Value = '  -0.22%  '

# --- Row 18: Bitcoin Cash ---
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '230.31'
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -4.66%  '

# --- Row 19: Shiba Inu ---
$ws.Cells.Item(19, 4).Value = '0.0₃0725'
$ws.Cells.Item(19, 5).Value = '  -0.69%  '

# --- Row 20: Chainlink ---
$ws.Cells.Item(20, 5).Value = '  -0.58%  '

# --- Row 22: Uniswap ---
$ws.Cells.Item(22, 5).Value = '  -3.60%  '

# --- Row 23: Avalanche ---
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '9.76'
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +4.00%  '

# --- Row 24: Toncoin ---
$ws.Cells.Item(24, 5).Value = '  -2.06%  '

# --- Row 25: Monero ---
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '148.65'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  +1.77%  '

# --- Row 26: Cosmos ---
$ws.Cells.Item(26, 5).Value = '  -2.62%  '

# --- Row 27: Stellar ---
$ws.Cells.Item(27, 5).Value = '  +1.33%  '

# --- Row 28: BinanceUSD ---
$ws.Cells.Item(28, 5).Value = '  +0.23%  '

# --- Row 29: EthereumClassic ---
$ws.Cells.Item(29, 5).Value = '  -4.34%  '

# --- Row 30: PancakeSwap ---
$ws.Cells.Item(30, 5).Value = '  -2.56%  '

# --- Row 31: Hedera ---
$ws.Cells.Item(31, 5).Value = '  -3.08%  '

# --- Row 32: Filecoin ---
$ws.Cells.Item(32, 5).Value = '  -0.72%  '

# --- Row 33: Internet Computer (DFINITY) ---
$ws.Cells.Item(33, 5).Value = '  +2.44%  '

# --- Row 34: Maker ---
$ws.Cells.Item(34, 4).Value = '1.426.50'
$ws.Cells.Item(34, 5).Value = '  -2.05%  '

# --- Row 35: LidoDAOToken ---
$ws.Cells.Item(35, 5).Value = '  +1.48%  '

# --- Row 36: HuobiToken ---
$ws.Cells.Item(36, 5).Value = '  -0.18%  '

# --- Row 37: ImmutableX ---
$ws.Cells.Item(37, 5).Value = '  -0.85%  '

# --- Row 38: ARBITRUM ---
$ws.Cells.Item(38, 5).Value = '  -4.42%  '

# --- Row 39: VeChain ---
$ws.Cells.Item(39, 5).Value = '  -3.08%  '

# --- Row 40: WEMIXToken ---
$ws.Cells.Item(40, 5).Value = '  +0.69%  '

# --- Row 41: PaxDollar ---
$ws.Cells.Item(41, 5).Value = '  +0.08%  '

# --- Row 42: TrustWalletToken ---
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.817'
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +3.02%  '

# --- Row 43: FraxShare ---
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '5.55'
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +2.47%  '

# --- Row 44: MXToken ---
$ws.Cells.Item(44, 5).Value = '  +1.27%  '

# --- Row 45: Aave ---
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '65.19'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -6.61%  '

# --- Row 46: RocketPoolETH ---
$ws.Cells.Item(46, 4).Value = '1.790.59'
$ws.Cells.Item(46, 5).Value = '  -0.58%  '

# --- Row 47: RenderToken ---
$ws.Cells.Item(47, 5).Value = '  -1.56%  '

# --- Row 48: Quant ---
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '88.13'
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.68%  '

# --- Row 49: BabyDogeCoin ---
$ws.Cells.Item(49, 5).Value = '  +0.74%  '

# --- Row 50: EnergySwap/Algorand swap ---
$ws.Cells.Item(50, 2).Value = 'Algorand'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0995'
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -3.03%  '

# --- Row 51: EnergySwap/Algorand swap ---
$ws.Cells.Item(51, 2).Value = 'EnergySwap'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '7.79'
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -0.58%  '

Write-Output "Updated cryptos list with refreshed prices and volumes."
